# Update "Date Dropped" column (G) for candidates who have since dropped out.
# Cells previously held the placeholder text "-"; now they get an actual
# drop date, formatted like the other already-dated rows (numFmtId 14,
# m/d/yyyy) by copying that formatting from an existing dated cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormatSource = $ws.Range("G3")

function Set-DropDate($cellRef, $serial) {
    $target = $ws.Range($cellRef)
    $dateFormatSource.Copy()
    $target.PasteSpecial(-4122)  # xlPasteFormats
    $target.Value = $serial
}

Set-DropDate "G2" 43872
Set-DropDate "G8" 43891
Set-DropDate "G17" 43892
Set-DropDate "G23" 43872

# Reflect the new active-cell selection recorded in the saved workbook.
$ws.Range("I18").Select()
